# Apply the reviewer's follow-up edit to the supplementary-information
# document:
#   1. Drop the stray "_GoBack" bookmark that sat at the very start of the
#      document (an artifact Word leaves behind at the last edit point).
#   2. Add a new tracked insertion ("the dominant ") by a new reviewer,
#      "German Uritskiy", right before the already-tracked deletion of
#      "major differentially abundant phyla" in the first OTU-relative-
#      abundance sentence (the Fig. S2/S4 caption).
#   3. Re-create the "_GoBack" bookmark at the new last-edit location,
#      splitting the "Error bars represent standard deviation; significance "
#      run right after "signific".
#
# The document's revisions keep their original author/date; only the new
# insertion (and the new bookmark position) is attributed to the new
# reviewer at the current time, matching how Word stamps freshly-typed
# content.

$d = $word.ActiveDocument

# Make sure the new edits are tracked, and attributed to the new reviewer
# mentioned in word/people.xml for this change.
$d.TrackRevisions = $true
$word.UserName = "German Uritskiy"

# --- 1. Remove the old "_GoBack" bookmark -----------------------------
# "_GoBack" is a hidden bookmark (leading underscore) so it will not show
# up in the normal Bookmarks enumeration, but it can still be addressed
# directly by name.
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
    # already absent - nothing to do
}

# --- 2. Insert "the dominant " before the tracked deletion of ---------
#        "major differentially abundant phyla" (first occurrence, the one
#        that reads "...(A-D) relative abundance of major differentially
#        abundant phyla").
$rng = $d.Content
$found = $rng.Find.Execute(
    "(A-D) relative abundance of", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    # Inserting exactly at the boundary between this run and the existing
    # <w:del> that follows it would be swallowed into/after the deletion,
    # so anchor one character earlier (the space before "of" ends) and
    # insert " the dominant" - the net visible text is identical
    # ("...relative abundance of the dominant phyla...") while the new
    # <w:ins> lands ahead of the old <w:del> as intended.
    $insertPos = $rng.End
    $target = $d.Range($insertPos, $insertPos)
    $target.InsertAfter(" the dominant")
}

# --- 3. Re-create "_GoBack" at the new last-edited location -----------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "Error bars represent standard deviation; signific", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $pos2 = $rng2.End
    $d.Bookmarks.Add("_GoBack", $d.Range($pos2, $pos2))
}
